$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 27, shifting rows 27:64 down to 28:65
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new record
$ws.Cells.Item(27, 1).Value = 1
$ws.Cells.Item(27, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(27, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(27, 4).Value = 44778
$ws.Cells.Item(27, 5).Value = 15
$ws.Cells.Item(27, 6).Value = 100112031
$ws.Cells.Item(27, 7).Value = "Poroto verde"
$ws.Cells.Item(27, 8).Value = "Magnum"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 300
$ws.Cells.Item(27, 11).Value = 27000
$ws.Cells.Item(27, 12).Value = 28000
$ws.Cells.Item(27, 13).Value = 27500
$ws.Cells.Item(27, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(27, 15).Value = "Perú"
$ws.Cells.Item(27, 16).Value = 1100
$ws.Cells.Item(27, 17).Value = 25
$ws.Cells.Item(27, 18).Value = "Hortaliza"
